# Updated cryptos list on Sat Jun 24 15:23:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "30.633.70" or "1.003" which Excel would
# otherwise auto-detect and coerce into a numeric cell (losing the original
# text formatting / introducing floating point artifacts). Force these
# cells to be treated as Text so the written value is preserved verbatim,
# matching the source workbook where every cell is an inline/shared string.
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "30.633.70"
$ws.Cells.Item(2, 5).Value = "  +0.08%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "1.889.60"
$ws.Cells.Item(3, 5).Value = "  -0.19%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).Value = "1.003"
$ws.Cells.Item(4, 5).Value = "  +0.61%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "240.67"
$ws.Cells.Item(5, 5).Value = "  -1.62%  "

# Row 6 - USDC
$ws.Cells.Item(6, 4).Value = "1.002"
$ws.Cells.Item(6, 5).Value = "  +0.54%  "

# Row 7 - XRP
$ws.Cells.Item(7, 4).Value = "0.4886"
$ws.Cells.Item(7, 5).Value = "  -0.81%  "

# Row 8 - Cardano
$ws.Cells.Item(8, 4).Value = "0.2914"
$ws.Cells.Item(8, 5).Value = "  -0.70%  "

# Row 9 - Dogecoin
$ws.Cells.Item(9, 4).Value = "0.06696"
$ws.Cells.Item(9, 5).Value = "  +0.67%  "

# Row 10 - WrappedEther
$ws.Cells.Item(10, 4).Value = "1.895.32"
$ws.Cells.Item(10, 5).Value = "  +0.76%  "

# Row 11 - Solana
$ws.Cells.Item(11, 4).Value = "17.00"
$ws.Cells.Item(11, 5).Value = "  +1.79%  "

# Row 12 - TRON
$ws.Cells.Item(12, 4).Value = "0.07238"
$ws.Cells.Item(12, 5).Value = "  +0.17%  "

# Row 13 - Litecoin
$ws.Cells.Item(13, 4).Value = "89.64"
$ws.Cells.Item(13, 5).Value = "  +3.19%  "

# Row 14 - was Polygon, now Polkadot
$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).Value = "5.016"
$ws.Cells.Item(14, 5).Value = "  -0.46%  "

# Row 15 - was Polkadot, now Polygon
$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(15, 4).Value = "0.6704"
$ws.Cells.Item(15, 5).Value = "  -0.68%  "

# Row 16 - WrappedBTC
$ws.Cells.Item(16, 4).Value = "30.625.77"
$ws.Cells.Item(16, 5).Value = "  +0.64%  "

# Row 17 - ShibaInu
$ws.Cells.Item(17, 4).Value = "0.000007936"
$ws.Cells.Item(17, 5).Value = "  +0.56%  "

# Row 18 - Dai
$ws.Cells.Item(18, 4).Value = "1.002"
$ws.Cells.Item(18, 5).Value = "  +0.25%  "

# Row 19 - Avalanche
$ws.Cells.Item(19, 4).Value = "13.09"
$ws.Cells.Item(19, 5).Value = "  +1.63%  "

# Row 20 - WrappedliquidstakedEther2.0
$ws.Cells.Item(20, 4).Value = "2.139.86"
$ws.Cells.Item(20, 5).Value = "  +0.82%  "

# Row 21 - BinanceUSD
$ws.Cells.Item(21, 4).Value = "1.002"
$ws.Cells.Item(21, 5).Value = "  +0.51%  "

# Row 22 - Uniswap
$ws.Cells.Item(22, 4).Value = "4.783"
$ws.Cells.Item(22, 5).Value = "  -0.27%  "

# Row 23 - BitcoinCash
$ws.Cells.Item(23, 4).Value = "189.62"
$ws.Cells.Item(23, 5).Value = "  +30.07%  "

# Row 24 - Chainlink
$ws.Cells.Item(24, 4).Value = "6.046"
$ws.Cells.Item(24, 5).Value = "  +1.76%  "

# Row 25 - Cosmos
$ws.Cells.Item(25, 4).Value = "9.342"
$ws.Cells.Item(25, 5).Value = "  +0.97%  "

# Row 26 - Monero
$ws.Cells.Item(26, 4).Value = "156.44"
$ws.Cells.Item(26, 5).Value = "  +2.85%  "

# Row 27 - EthereumClassic
$ws.Cells.Item(27, 4).Value = "18.60"
$ws.Cells.Item(27, 5).Value = "  +8.82%  "

# Row 28 - LidoDAOToken
$ws.Cells.Item(28, 4).Value = "1.879"
$ws.Cells.Item(28, 5).Value = "  -2.24%  "

# Row 29 - Toncoin
$ws.Cells.Item(29, 4).Value = "1.414"
$ws.Cells.Item(29, 5).Value = "  -1.05%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Cells.Item(30, 4).Value = "4.266"
$ws.Cells.Item(30, 5).Value = "  +0.28%  "

# Row 31 - Stellar
$ws.Cells.Item(31, 4).Value = "0.09042"
$ws.Cells.Item(31, 5).Value = "  +2.52%  "

# Row 32 - Filecoin
$ws.Cells.Item(32, 4).Value = "3.962"
$ws.Cells.Item(32, 5).Value = "  -1.24%  "

# Row 33 - Hedera
$ws.Cells.Item(33, 4).Value = "0.05295"
$ws.Cells.Item(33, 5).Value = "  +1.31%  "

# Row 34 - ImmutableX
$ws.Cells.Item(34, 4).Value = "0.7386"
$ws.Cells.Item(34, 5).Value = "  +1.03%  "

# Row 35 - ARBITRUM
$ws.Cells.Item(35, 4).Value = "1.092"
$ws.Cells.Item(35, 5).Value = "  -3.04%  "

# Row 36 - HuobiToken
$ws.Cells.Item(36, 4).Value = "2.767"
$ws.Cells.Item(36, 5).Value = "  +4.02%  "

# Row 37 - VeChain
$ws.Cells.Item(37, 4).Value = "0.01823"
$ws.Cells.Item(37, 5).Value = "  -1.16%  "

# Row 38 - MXToken
$ws.Cells.Item(38, 4).Value = "2.676"
$ws.Cells.Item(38, 5).Value = "  -1.09%  "

# Row 39 - TrustWalletToken
$ws.Cells.Item(39, 4).Value = "0.9247"
$ws.Cells.Item(39, 5).Value = "  -1.80%  "

# Row 40 - RenderToken
$ws.Cells.Item(40, 4).Value = "2.077"
$ws.Cells.Item(40, 5).Value = "  -5.10%  "

# Row 41 - TheSandbox
$ws.Cells.Item(41, 4).Value = "0.4386"
$ws.Cells.Item(41, 5).Value = "  +2.12%  "

# Row 42 - Quant
$ws.Cells.Item(42, 4).Value = "105.05"
$ws.Cells.Item(42, 5).Value = "  +0.82%  "

# Row 43 - PaxDollar (D unchanged)
$ws.Cells.Item(43, 5).Value = "  +0.26%  "

# Row 44 - FraxShare
$ws.Cells.Item(44, 4).Value = "5.700"
$ws.Cells.Item(44, 5).Value = "  -2.27%  "

# Row 45 - Algorand
$ws.Cells.Item(45, 4).Value = "0.1348"
$ws.Cells.Item(45, 5).Value = "  +3.28%  "

# Row 46 - Aptos
$ws.Cells.Item(46, 4).Value = "7.451"
$ws.Cells.Item(46, 5).Value = "  -1.65%  "

# Row 47 - Cronos
$ws.Cells.Item(47, 4).Value = "0.05862"
$ws.Cells.Item(47, 5).Value = "  +1.69%  "

# Row 48 - EnergySwap
$ws.Cells.Item(48, 4).Value = "8.811"
$ws.Cells.Item(48, 5).Value = "  +6.01%  "

# Row 49 - was Decentraland, now Elrond
$ws.Cells.Item(49, 2).Value = "Elrond"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(49, 4).Value = "33.80"
$ws.Cells.Item(49, 5).Value = "  +2.45%  "

# Row 50 - was Elrond, now Decentraland
$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(50, 4).Value = "0.3941"
$ws.Cells.Item(50, 5).Value = "  +3.43%  "

# Row 51 - NEARProtocol
$ws.Cells.Item(51, 4).Value = "1.419"
$ws.Cells.Item(51, 5).Value = "  +4.73%  "
